$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308:399 down to 309:400
$ws.Rows("308:308").Insert()

# Populate the newly inserted row 308 with the new weekly record
$ws.Range("A308").Value = 9
$ws.Range("B308").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C308").Value = "Metropolitana"
$ws.Range("D308").Value = 45215
$ws.Range("E308").Value = 13
$ws.Range("F308").Value = 100112001
$ws.Range("G308").Value = "Berenjena"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 124
$ws.Range("K308").Value = 9000
$ws.Range("L308").Value = 10000
$ws.Range("M308").Value = 9500
$ws.Range("N308").Value = "$/caja 50 unidades"
$ws.Range("O308").Value = "Región de Arica y Parinacota"
$ws.Range("P308").Value = 190
$ws.Range("Q308").Value = 50
$ws.Range("R308").Value = "Hortaliza"
